# RPA datasets push 2023-10-20
# Rewrites the data rows (2-26) of sheet1 to match the refreshed dataset export,
# inserting/reordering rows for 2023-10-10 (신성에스티) and 2023-09-25 (에이치엠씨제6호스팩),
# and re-sequencing the 두산로보틱스 / 유안타제11호스팩 rows that moved position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, F, G hold dates stored as plain text (e.g. "2023-09-21").
# Force text format first so COM does not silently coerce them to date serials,
# then restore the default "Normal" style once the values are in place so the
# on-disk cell styling matches the original (unstyled) data rows.
$dateColB = $ws.Range("B2:B26")
$dateColF = $ws.Range("F2:F26")
$dateColG = $ws.Range("G2:G26")
$dateColB.NumberFormat = "@"
$dateColF.NumberFormat = "@"
$dateColG.NumberFormat = "@"

# Row 2: 두산로보틱스
$ws.Cells.Item(2, 1).Value = 'CS'
$ws.Cells.Item(2, 2).Value = '2023-09-21'
$ws.Cells.Item(2, 3).Value = '두산로보틱스'
$ws.Cells.Item(2, 4).Value = '한국, 미래'
$ws.Cells.Item(2, 5).Value = '한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)'
$ws.Cells.Item(2, 6).Value = '2023-09-26'
$ws.Cells.Item(2, 7).Value = '2023-10-05'
$ws.Cells.Item(2, 8).Value = 42120
$ws.Cells.Item(2, 9).Value = 16200000
$ws.Cells.Item(2, 10).Value = 26000
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 10

# Row 3: 한싹
$ws.Cells.Item(3, 1).Value = 'KB'
$ws.Cells.Item(3, 2).Value = '2023-09-19'
$ws.Cells.Item(3, 3).Value = '한싹'
$ws.Cells.Item(3, 4).Value = 'KB'
$ws.Cells.Item(3, 5).Value = 'KB'
$ws.Cells.Item(3, 6).Value = '2023-09-22'
$ws.Cells.Item(3, 7).Value = '2023-10-04'
$ws.Cells.Item(3, 8).Value = 18750
$ws.Cells.Item(3, 9).Value = 1500000
$ws.Cells.Item(3, 10).Value = 12500
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 100

# Row 4: 두산로보틱스
$ws.Cells.Item(4, 1).Value = 'KB'
$ws.Cells.Item(4, 2).Value = '2023-09-21'
$ws.Cells.Item(4, 3).Value = '두산로보틱스'
$ws.Cells.Item(4, 4).Value = '한국, 미래'
$ws.Cells.Item(4, 5).Value = '한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)'
$ws.Cells.Item(4, 6).Value = '2023-09-26'
$ws.Cells.Item(4, 7).Value = '2023-10-05'
$ws.Cells.Item(4, 8).Value = 42120
$ws.Cells.Item(4, 9).Value = 16200000
$ws.Cells.Item(4, 10).Value = 26000
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 10

# Row 5: 두산로보틱스
$ws.Cells.Item(5, 1).Value = 'NH'
$ws.Cells.Item(5, 2).Value = '2023-09-21'
$ws.Cells.Item(5, 3).Value = '두산로보틱스'
$ws.Cells.Item(5, 4).Value = '한국, 미래'
$ws.Cells.Item(5, 5).Value = '한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)'
$ws.Cells.Item(5, 6).Value = '2023-09-26'
$ws.Cells.Item(5, 7).Value = '2023-10-05'
$ws.Cells.Item(5, 8).Value = 42120
$ws.Cells.Item(5, 9).Value = 16200000
$ws.Cells.Item(5, 10).Value = 26000
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 10

# Row 6: 대신밸런스제15호스팩
$ws.Cells.Item(6, 1).Value = '대신'
$ws.Cells.Item(6, 2).Value = '2023-08-21'
$ws.Cells.Item(6, 3).Value = '대신밸런스제15호스팩'
$ws.Cells.Item(6, 4).Value = '대신'
$ws.Cells.Item(6, 5).Value = '대신'
$ws.Cells.Item(6, 6).Value = '2023-08-24'
$ws.Cells.Item(6, 7).Value = '2023-08-30'
$ws.Cells.Item(6, 8).Value = 13000
$ws.Cells.Item(6, 9).Value = 6500000
$ws.Cells.Item(6, 10).Value = 2000
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 100

# Row 7: 대신밸런스제16호스팩
$ws.Cells.Item(7, 1).Value = '대신'
$ws.Cells.Item(7, 2).Value = '2023-08-23'
$ws.Cells.Item(7, 3).Value = '대신밸런스제16호스팩'
$ws.Cells.Item(7, 4).Value = '대신'
$ws.Cells.Item(7, 5).Value = '대신'
$ws.Cells.Item(7, 6).Value = '2023-08-28'
$ws.Cells.Item(7, 7).Value = '2023-09-04'
$ws.Cells.Item(7, 8).Value = 13000
$ws.Cells.Item(7, 9).Value = 6500000
$ws.Cells.Item(7, 10).Value = 2000
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 100

# Row 8: 밀리의서재
$ws.Cells.Item(8, 1).Value = '미래'
$ws.Cells.Item(8, 2).Value = '2023-09-18'
$ws.Cells.Item(8, 3).Value = '밀리의서재'
$ws.Cells.Item(8, 4).Value = '미래'
$ws.Cells.Item(8, 5).Value = '미래'
$ws.Cells.Item(8, 6).Value = '2023-09-21'
$ws.Cells.Item(8, 7).Value = '2023-09-27'
$ws.Cells.Item(8, 8).Value = 34500
$ws.Cells.Item(8, 9).Value = 1500000
$ws.Cells.Item(8, 10).Value = 23000
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 100

# Row 9: 신성에스티
$ws.Cells.Item(9, 1).Value = '미래'
$ws.Cells.Item(9, 2).Value = '2023-10-10'
$ws.Cells.Item(9, 3).Value = '신성에스티'
$ws.Cells.Item(9, 4).Value = '미래'
$ws.Cells.Item(9, 5).Value = '미래'
$ws.Cells.Item(9, 6).Value = '2023-10-13'
$ws.Cells.Item(9, 7).Value = '2023-10-19'
$ws.Cells.Item(9, 8).Value = 52000
$ws.Cells.Item(9, 9).Value = 2000000
$ws.Cells.Item(9, 10).Value = 26000
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 100

# Row 10: 퓨릿
$ws.Cells.Item(10, 1).Value = '미래'
$ws.Cells.Item(10, 2).Value = '2023-10-05'
$ws.Cells.Item(10, 3).Value = '퓨릿'
$ws.Cells.Item(10, 4).Value = '미래'
$ws.Cells.Item(10, 5).Value = '미래'
$ws.Cells.Item(10, 6).Value = '2023-10-11'
$ws.Cells.Item(10, 7).Value = '2023-10-18'
$ws.Cells.Item(10, 8).Value = 44265.9
$ws.Cells.Item(10, 9).Value = 4137000
$ws.Cells.Item(10, 10).Value = 10700
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 100

# Row 11: 두산로보틱스
$ws.Cells.Item(11, 1).Value = '미래'
$ws.Cells.Item(11, 2).Value = '2023-09-21'
$ws.Cells.Item(11, 3).Value = '두산로보틱스'
$ws.Cells.Item(11, 4).Value = '한국, 미래'
$ws.Cells.Item(11, 5).Value = '한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)'
$ws.Cells.Item(11, 6).Value = '2023-09-26'
$ws.Cells.Item(11, 7).Value = '2023-10-05'
$ws.Cells.Item(11, 8).Value = 126360
$ws.Cells.Item(11, 9).Value = 16200000
$ws.Cells.Item(11, 10).Value = 26000
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 30

# Row 12: 레뷰코퍼레이션
$ws.Cells.Item(12, 1).Value = '삼성'
$ws.Cells.Item(12, 2).Value = '2023-09-19'
$ws.Cells.Item(12, 3).Value = '레뷰코퍼레이션'
$ws.Cells.Item(12, 4).Value = '삼성'
$ws.Cells.Item(12, 5).Value = '삼성'
$ws.Cells.Item(12, 6).Value = '2023-09-22'
$ws.Cells.Item(12, 7).Value = '2023-10-06'
$ws.Cells.Item(12, 8).Value = 33600
$ws.Cells.Item(12, 9).Value = 2240000
$ws.Cells.Item(12, 10).Value = 15000
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 100

# Row 13: 상상인제4호스팩
$ws.Cells.Item(13, 1).Value = '상상인'
$ws.Cells.Item(13, 2).Value = '2023-09-04'
$ws.Cells.Item(13, 3).Value = '상상인제4호스팩'
$ws.Cells.Item(13, 4).Value = '상상인'
$ws.Cells.Item(13, 5).Value = '상상인'
$ws.Cells.Item(13, 6).Value = '2023-09-07'
$ws.Cells.Item(13, 7).Value = '2023-09-14'
$ws.Cells.Item(13, 8).Value = 9000
$ws.Cells.Item(13, 9).Value = 4500000
$ws.Cells.Item(13, 10).Value = 2000
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 100

# Row 14: 두산로보틱스
$ws.Cells.Item(14, 1).Value = '신영'
$ws.Cells.Item(14, 2).Value = '2023-09-21'
$ws.Cells.Item(14, 3).Value = '두산로보틱스'
$ws.Cells.Item(14, 4).Value = '한국, 미래'
$ws.Cells.Item(14, 5).Value = '한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)'
$ws.Cells.Item(14, 6).Value = '2023-09-26'
$ws.Cells.Item(14, 7).Value = '2023-10-05'
$ws.Cells.Item(14, 8).Value = 12636
$ws.Cells.Item(14, 9).Value = 16200000
$ws.Cells.Item(14, 10).Value = 26000
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 3

# Row 15: 인스웨이브시스템즈
$ws.Cells.Item(15, 1).Value = '신영'
$ws.Cells.Item(15, 2).Value = '2023-09-14'
$ws.Cells.Item(15, 3).Value = '인스웨이브시스템즈'
$ws.Cells.Item(15, 4).Value = '신영'
$ws.Cells.Item(15, 5).Value = '신영'
$ws.Cells.Item(15, 6).Value = '2023-09-19'
$ws.Cells.Item(15, 7).Value = '2023-09-25'
$ws.Cells.Item(15, 8).Value = 26400
$ws.Cells.Item(15, 9).Value = 1100000
$ws.Cells.Item(15, 10).Value = 24000
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 100

# Row 16: 신한제11호스팩
$ws.Cells.Item(16, 1).Value = '신한'
$ws.Cells.Item(16, 2).Value = '2023-09-19'
$ws.Cells.Item(16, 3).Value = '신한제11호스팩'
$ws.Cells.Item(16, 4).Value = '신한'
$ws.Cells.Item(16, 5).Value = '신한'
$ws.Cells.Item(16, 6).Value = '2023-09-22'
$ws.Cells.Item(16, 7).Value = '2023-10-04'
$ws.Cells.Item(16, 8).Value = 36000
$ws.Cells.Item(16, 9).Value = 18000000
$ws.Cells.Item(16, 10).Value = 2000
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 100

# Row 17: 두산로보틱스
$ws.Cells.Item(17, 1).Value = '유비에스'
$ws.Cells.Item(17, 2).Value = '2023-09-21'
$ws.Cells.Item(17, 3).Value = '두산로보틱스'
$ws.Cells.Item(17, 4).Value = '한국, 미래'
$ws.Cells.Item(17, 5).Value = '한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)'
$ws.Cells.Item(17, 6).Value = '2023-09-26'
$ws.Cells.Item(17, 7).Value = '2023-10-05'
$ws.Cells.Item(17, 8).Value = 4212
$ws.Cells.Item(17, 9).Value = 16200000
$ws.Cells.Item(17, 10).Value = 26000
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 1

# Row 18: 유안타제11호스팩
$ws.Cells.Item(18, 1).Value = '유안타'
$ws.Cells.Item(18, 2).Value = '2023-08-22'
$ws.Cells.Item(18, 3).Value = '유안타제11호스팩'
$ws.Cells.Item(18, 4).Value = '유안타'
$ws.Cells.Item(18, 5).Value = '유안타'
$ws.Cells.Item(18, 6).Value = '2023-08-25'
$ws.Cells.Item(18, 7).Value = '2023-09-01'
$ws.Cells.Item(18, 8).Value = 10000
$ws.Cells.Item(18, 9).Value = 5000000
$ws.Cells.Item(18, 10).Value = 2000
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 12).Value = 100

# Row 19: 아이엠티
$ws.Cells.Item(19, 1).Value = '유안타'
$ws.Cells.Item(19, 2).Value = '2023-09-18'
$ws.Cells.Item(19, 3).Value = '아이엠티'
$ws.Cells.Item(19, 4).Value = '유안타'
$ws.Cells.Item(19, 5).Value = '유안타, 유진'
$ws.Cells.Item(19, 6).Value = '2023-09-21'
$ws.Cells.Item(19, 7).Value = '2023-10-10'
$ws.Cells.Item(19, 8).Value = 15484
$ws.Cells.Item(19, 9).Value = 1580000
$ws.Cells.Item(19, 10).Value = 14000
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 12).Value = 70

# Row 20: 아이엠티
$ws.Cells.Item(20, 1).Value = '유진'
$ws.Cells.Item(20, 2).Value = '2023-09-18'
$ws.Cells.Item(20, 3).Value = '아이엠티'
$ws.Cells.Item(20, 4).Value = '유안타'
$ws.Cells.Item(20, 5).Value = '유안타, 유진'
$ws.Cells.Item(20, 6).Value = '2023-09-21'
$ws.Cells.Item(20, 7).Value = '2023-10-10'
$ws.Cells.Item(20, 8).Value = 6636
$ws.Cells.Item(20, 9).Value = 1580000
$ws.Cells.Item(20, 10).Value = 14000
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = 30

# Row 21: 두산로보틱스
$ws.Cells.Item(21, 1).Value = '키움'
$ws.Cells.Item(21, 2).Value = '2023-09-21'
$ws.Cells.Item(21, 3).Value = '두산로보틱스'
$ws.Cells.Item(21, 4).Value = '한국, 미래'
$ws.Cells.Item(21, 5).Value = '한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)'
$ws.Cells.Item(21, 6).Value = '2023-09-26'
$ws.Cells.Item(21, 7).Value = '2023-10-05'
$ws.Cells.Item(21, 8).Value = 12636
$ws.Cells.Item(21, 9).Value = 16200000
$ws.Cells.Item(21, 10).Value = 26000
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = 3

# Row 22: 두산로보틱스
$ws.Cells.Item(22, 1).Value = '하나'
$ws.Cells.Item(22, 2).Value = '2023-09-21'
$ws.Cells.Item(22, 3).Value = '두산로보틱스'
$ws.Cells.Item(22, 4).Value = '한국, 미래'
$ws.Cells.Item(22, 5).Value = '한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)'
$ws.Cells.Item(22, 6).Value = '2023-09-26'
$ws.Cells.Item(22, 7).Value = '2023-10-05'
$ws.Cells.Item(22, 8).Value = 12636
$ws.Cells.Item(22, 9).Value = 16200000
$ws.Cells.Item(22, 10).Value = 26000
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 3

# Row 23: 두산로보틱스
$ws.Cells.Item(23, 1).Value = '한국'
$ws.Cells.Item(23, 2).Value = '2023-09-21'
$ws.Cells.Item(23, 3).Value = '두산로보틱스'
$ws.Cells.Item(23, 4).Value = '한국, 미래'
$ws.Cells.Item(23, 5).Value = '한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)'
$ws.Cells.Item(23, 6).Value = '2023-09-26'
$ws.Cells.Item(23, 7).Value = '2023-10-05'
$ws.Cells.Item(23, 8).Value = 126360
$ws.Cells.Item(23, 9).Value = 16200000
$ws.Cells.Item(23, 10).Value = 26000
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 30

# Row 24: 한국제12호스팩
$ws.Cells.Item(24, 1).Value = '한국'
$ws.Cells.Item(24, 2).Value = '2023-08-21'
$ws.Cells.Item(24, 3).Value = '한국제12호스팩'
$ws.Cells.Item(24, 4).Value = '한국'
$ws.Cells.Item(24, 5).Value = '한국'
$ws.Cells.Item(24, 6).Value = '2023-08-24'
$ws.Cells.Item(24, 7).Value = '2023-08-30'
$ws.Cells.Item(24, 8).Value = 8000
$ws.Cells.Item(24, 9).Value = 4000000
$ws.Cells.Item(24, 10).Value = 2000
$ws.Cells.Item(24, 11).Value = 0
$ws.Cells.Item(24, 12).Value = 100

# Row 25: 한화플러스제4호스팩
$ws.Cells.Item(25, 1).Value = '한화'
$ws.Cells.Item(25, 2).Value = '2023-08-29'
$ws.Cells.Item(25, 3).Value = '한화플러스제4호스팩'
$ws.Cells.Item(25, 4).Value = '한화'
$ws.Cells.Item(25, 5).Value = '한화'
$ws.Cells.Item(25, 6).Value = '2023-09-01'
$ws.Cells.Item(25, 7).Value = '2023-09-07'
$ws.Cells.Item(25, 8).Value = 9500
$ws.Cells.Item(25, 9).Value = 4750000
$ws.Cells.Item(25, 10).Value = 2000
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 100

# Row 26: 에이치엠씨제6호스팩
$ws.Cells.Item(26, 1).Value = '현대차'
$ws.Cells.Item(26, 2).Value = '2023-09-25'
$ws.Cells.Item(26, 3).Value = '에이치엠씨제6호스팩'
$ws.Cells.Item(26, 4).Value = '현대차'
$ws.Cells.Item(26, 5).Value = '현대차'
$ws.Cells.Item(26, 6).Value = '2023-10-04'
$ws.Cells.Item(26, 7).Value = '2023-10-13'
$ws.Cells.Item(26, 8).Value = 8000
$ws.Cells.Item(26, 9).Value = 4000000
$ws.Cells.Item(26, 10).Value = 2000
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = 100

# Restore default styling on the date columns now that the text values are set.
$dateRange.Style = "Normal"
